$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText {
    param($sheet, [string]$cellAddr, [string]$text)
    $scratch = $sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($cellAddr).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
    $scratch.Clear()
}

$ws.Range('E2').Value = '2026-02-23 22:18:45'
$ws.Range('O2').Value = '5.7 °C'
$ws.Range('E3').Value = '2026-02-23 22:18:48'
$ws.Range('E4').Value = '2026-02-23 22:18:51'
Set-PercentText $ws 'H4' '70%'
$ws.Range('O4').Value = '11.9 °C'
$ws.Range('E5').Value = '2026-02-23 22:18:53'
$ws.Range('E6').Value = '2026-02-23 22:18:56'
$ws.Range('O6').Value = '13.9 °C'
$ws.Range('E7').Value = '2026-02-23 22:18:58'
$ws.Range('O7').Value = '14.1 °C'
$ws.Range('E8').Value = '2026-02-23 22:19:00'
$ws.Range('E9').Value = '2026-02-23 22:19:03'
$ws.Range('E10').Value = '2026-02-23 22:19:06'
Set-PercentText $ws 'H10' '77%'
$ws.Range('O10').Value = '10.6 °C'
$ws.Range('E11').Value = '2026-02-23 22:19:08'
$ws.Range('O11').Value = '8.7 °C'
$ws.Range('E12').Value = '2026-02-23 22:19:11'
Set-PercentText $ws 'H12' '86%'
$ws.Range('E13').Value = '2026-02-23 22:19:13'
$ws.Range('E14').Value = '2026-02-23 22:19:16'
Set-PercentText $ws 'H14' '77%'
$ws.Range('O14').Value = '12.3 °C'
$ws.Range('E15').Value = '2026-02-23 22:19:19'
$ws.Range('O15').Value = '12.3 °C'
$ws.Range('E16').Value = '2026-02-23 22:19:21'
Set-PercentText $ws 'H16' '23%'
$ws.Range('E17').Value = '2026-02-23 22:19:23'
$ws.Range('E18').Value = '2026-02-23 22:19:26'
$ws.Range('O18').Value = '10.8 °C'
$ws.Range('E19').Value = '2026-02-23 22:19:28'
Set-PercentText $ws 'H19' '48%'
$ws.Range('E20').Value = '2026-02-23 22:19:31'
$ws.Range('E21').Value = '2026-02-23 22:19:34'
$ws.Range('K21').Value = '16.2 MJ/m2'
$ws.Range('E22').Value = '2026-02-23 22:19:36'
$ws.Range('E23').Value = '2026-02-23 22:19:39'
$ws.Range('E24').Value = '2026-02-23 22:19:42'
$ws.Range('E25').Value = '2026-02-23 22:19:45'
$ws.Range('E26').Value = '2026-02-23 22:19:47'
$ws.Range('J26').Value = '1024.0 hPa'
$ws.Range('O26').Value = '9.9 °C'
$ws.Range('E27').Value = '2026-02-23 22:19:50'
$ws.Range('E28').Value = '2026-02-23 22:19:52'
Set-PercentText $ws 'H28' '69%'
$ws.Range('O28').Value = '10.9 °C'
$ws.Range('E29').Value = '2026-02-23 22:19:55'
$ws.Range('O29').Value = '10.6 °C'
$ws.Range('E30').Value = '2026-02-23 22:19:58'
Set-PercentText $ws 'H30' '71%'
$ws.Range('O30').Value = '12.8 °C'
$ws.Range('E31').Value = '2026-02-23 22:20:00'
Set-PercentText $ws 'H31' '45%'
$ws.Range('E32').Value = '2026-02-23 22:20:03'
$ws.Range('N32').Value = '-1.2 °C 21:41 TU'
$ws.Range('O32').Value = '7.2 °C'
$ws.Range('E33').Value = '2026-02-23 22:20:05'
Set-PercentText $ws 'H33' '46%'
$ws.Range('J33').Value = '1025.4 hPa'
$ws.Range('O33').Value = '8.6 °C'
$ws.Range('E34').Value = '2026-02-23 22:20:08'
$ws.Range('E35').Value = '2026-02-23 22:20:11'
$ws.Range('J35').Value = '1025.2 hPa'
$ws.Range('O35').Value = '12.0 °C'
$ws.Range('E36').Value = '2026-02-23 22:20:13'
$ws.Range('J36').Value = '1024.5 hPa'
$ws.Range('O36').Value = '12.9 °C'
$ws.Range('E37').Value = '2026-02-23 22:20:16'
Set-PercentText $ws 'H37' '68%'
$ws.Range('O37').Value = '8.9 °C'
$ws.Range('E38').Value = '2026-02-23 22:20:19'
$ws.Range('O38').Value = '12.1 °C'
$ws.Range('E39').Value = '2026-02-23 22:20:21'
$ws.Range('E40').Value = '2026-02-23 22:20:24'
Set-PercentText $ws 'H40' '62%'
$ws.Range('J40').Value = '1026.3 hPa'
$ws.Range('O40').Value = '8.6 °C'
$ws.Range('E41').Value = '2026-02-23 22:20:26'
Set-PercentText $ws 'H41' '74%'
$ws.Range('O41').Value = '11.8 °C'
$ws.Range('E42').Value = '2026-02-23 22:20:29'
Set-PercentText $ws 'H42' '80%'
$ws.Range('O42').Value = '11.7 °C'
$ws.Range('E43').Value = '2026-02-23 22:20:31'
$ws.Range('O43').Value = '10.4 °C'
$ws.Range('E44').Value = '2026-02-23 22:20:34'
$ws.Range('E45').Value = '2026-02-23 22:20:37'
$ws.Range('E46').Value = '2026-02-23 22:20:39'
$ws.Range('J46').Value = '1025.8 hPa'
$ws.Range('O46').Value = '10.1 °C'

$wb.Save()
